# Auto-applies numeric cell updates to the Ridill_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 66670830
$ws.Cells.Item(100, 9).Value = 2500
$ws.Cells.Item(100, 10).Value = 100005000
$ws.Cells.Item(100, 11).Value = 2500
$ws.Cells.Item(100, 12).Value = 100005000
$ws.Cells.Item(100, 13).Value = -1959
$ws.Cells.Item(100, 14).Value = -100006082
$ws.Cells.Item(113, 8).Value = 14287714
$ws.Cells.Item(113, 9).Value = 25001500
$ws.Cells.Item(113, 11).Value = 25001500
$ws.Cells.Item(113, 13).Value = -24998246
$ws.Cells.Item(141, 8).Value = 7248.75
$ws.Cells.Item(141, 9).Value = 2331.6667
$ws.Cells.Item(141, 11).Value = 6995.000100000001
$ws.Cells.Item(141, 13).Value = -1815.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 44851.375
$ws.Cells.Item(2, 9).Value = 51001.57
$ws.Cells.Item(2, 10).Value = 1800
$ws.Cells.Item(2, 11).Value = 51001.57
$ws.Cells.Item(2, 12).Value = 1800
$ws.Cells.Item(2, 13).Value = -50888.57
$ws.Cells.Item(2, 14).Value = -2026
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 40000
$ws.Cells.Item(109, 10).Value = 40000
$ws.Cells.Item(109, 12).Value = 40000
$ws.Cells.Item(109, 14).Value = -42774
$ws.Cells.Item(110, 8).Value = 831.6
$ws.Cells.Item(110, 9).Value = 547.9
$ws.Cells.Item(110, 10).Value = 1399
$ws.Cells.Item(110, 11).Value = 547.9
$ws.Cells.Item(110, 12).Value = 1399
$ws.Cells.Item(110, 13).Value = 1497.1
$ws.Cells.Item(110, 14).Value = -5489
$ws.Cells.Item(111, 8).Value = 42000
$ws.Cells.Item(111, 10).Value = 42000
$ws.Cells.Item(111, 12).Value = 42000
$ws.Cells.Item(111, 14).Value = -50180
$ws.Cells.Item(112, 8).Value = 42000
$ws.Cells.Item(112, 10).Value = 42000
$ws.Cells.Item(112, 12).Value = 42000
$ws.Cells.Item(112, 14).Value = -44954
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()
$ws.Cells.Item(116, 8).Value = 44851.375
$ws.Cells.Item(116, 9).Value = 51001.57
$ws.Cells.Item(116, 10).Value = 1800
$ws.Cells.Item(116, 11).Value = 51001.57
$ws.Cells.Item(116, 12).Value = 1800
$ws.Cells.Item(116, 13).Value = -48707.57
$ws.Cells.Item(116, 14).Value = -6388
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()
$ws.Cells.Item(118, 8).Value = 34306.75
$ws.Cells.Item(118, 10).Value = 34306.75
$ws.Cells.Item(118, 12).Value = 34306.75
$ws.Cells.Item(118, 14).Value = -37620.75
$ws.Cells.Item(119, 8).Value = 20383
$ws.Cells.Item(119, 10).Value = 20383
$ws.Cells.Item(119, 12).Value = 20383
$ws.Cells.Item(119, 14).Value = -30059
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 16209511
$ws.Cells.Item(132, 9).Value = 22230120
$ws.Cells.Item(132, 10).Value = 6175163
$ws.Cells.Item(132, 11).Value = 66690360
$ws.Cells.Item(132, 12).Value = 18525489
$ws.Cells.Item(132, 13).Value = -66687830
$ws.Cells.Item(132, 14).Value = -18530549

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 44851.375
$ws.Cells.Item(3, 9).Value = 51001.57
$ws.Cells.Item(3, 10).Value = 1800
$ws.Cells.Item(3, 11).Value = 51001.57
$ws.Cells.Item(3, 12).Value = 1800
$ws.Cells.Item(3, 13).Value = -50887.57
$ws.Cells.Item(3, 14).Value = -2028
$ws.Cells.Item(52, 8).Value = 31168.572
$ws.Cells.Item(52, 10).Value = 31168.572
$ws.Cells.Item(52, 12).Value = 31168.572
$ws.Cells.Item(52, 14).Value = -31694.572
$ws.Cells.Item(99, 8).Value = 1749.875
$ws.Cells.Item(99, 9).Value = 1333.1666
$ws.Cells.Item(99, 10).Value = 3000
$ws.Cells.Item(99, 11).Value = 1333.1666
$ws.Cells.Item(99, 12).Value = 3000
$ws.Cells.Item(99, 13).Value = 164.8334
$ws.Cells.Item(99, 14).Value = -5996
$ws.Cells.Item(107, 8).Value = 1027.75
$ws.Cells.Item(107, 9).Value = 1037
$ws.Cells.Item(107, 11).Value = 1037
$ws.Cells.Item(107, 13).Value = 883
$ws.Cells.Item(109, 8).Value = 39800
$ws.Cells.Item(109, 10).Value = 39800
$ws.Cells.Item(109, 12).Value = 39800
$ws.Cells.Item(109, 14).Value = -42574
$ws.Cells.Item(110, 8).Value = 42000
$ws.Cells.Item(110, 10).Value = 42000
$ws.Cells.Item(110, 12).Value = 42000
$ws.Cells.Item(110, 14).Value = -50180
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 30684
$ws.Cells.Item(115, 10).Value = 30684
$ws.Cells.Item(115, 12).Value = 30684
$ws.Cells.Item(115, 14).Value = -33818
$ws.Cells.Item(116, 8).Value = 38000
$ws.Cells.Item(116, 10).Value = 38000
$ws.Cells.Item(116, 12).Value = 38000
$ws.Cells.Item(116, 14).Value = -47178
$ws.Cells.Item(119, 8).Value = 35000
$ws.Cells.Item(119, 10).Value = 35000
$ws.Cells.Item(119, 12).Value = 35000
$ws.Cells.Item(119, 14).Value = -44676
$ws.Cells.Item(121, 8).Value = 31168.572
$ws.Cells.Item(121, 10).Value = 31168.572
$ws.Cells.Item(121, 12).Value = 31168.572
$ws.Cells.Item(121, 14).Value = -34662.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1022.2
$ws.Cells.Item(16, 9).Value = 1022.2
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1022.2
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -735.2
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(110, 8).Value = 35000
$ws.Cells.Item(110, 10).Value = 35000
$ws.Cells.Item(110, 12).Value = 35000
$ws.Cells.Item(110, 14).Value = -43180
$ws.Cells.Item(111, 8).Value = 42000
$ws.Cells.Item(111, 10).Value = 42000
$ws.Cells.Item(111, 12).Value = 42000
$ws.Cells.Item(111, 14).Value = -50180
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 1022.2
$ws.Cells.Item(113, 9).Value = 1022.2
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1022.2
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = 1147.8
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(115, 8).Value = 30001
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 30001
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 12).Value = 30001
$ws.Cells.Item(115, 13).ClearContents()
$ws.Cells.Item(115, 14).Value = -32351
$ws.Cells.Item(116, 8).Value = 22124.25
$ws.Cells.Item(116, 10).Value = 22124.25
$ws.Cells.Item(116, 12).Value = 22124.25
$ws.Cells.Item(116, 14).Value = -31302.25
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(119, 8).Value = 43120.332
$ws.Cells.Item(119, 10).Value = 43120.332
$ws.Cells.Item(119, 12).Value = 43120.332
$ws.Cells.Item(119, 14).Value = -52796.332
$ws.Cells.Item(120, 8).Value = 71674
$ws.Cells.Item(120, 9).Value = 200296
$ws.Cells.Item(120, 10).Value = 28800
$ws.Cells.Item(120, 11).Value = 200296
$ws.Cells.Item(120, 12).Value = 28800
$ws.Cells.Item(120, 13).Value = -196667
$ws.Cells.Item(120, 14).Value = -36058
$ws.Cells.Item(121, 8).Value = 19975
$ws.Cells.Item(121, 10).Value = 19975
$ws.Cells.Item(121, 12).Value = 19975
$ws.Cells.Item(121, 14).Value = -22595

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 222.90909
$ws.Cells.Item(107, 9).Value = 137.42857
$ws.Cells.Item(107, 10).Value = 372.5
$ws.Cells.Item(107, 11).Value = 137.42857
$ws.Cells.Item(107, 12).Value = 372.5
$ws.Cells.Item(107, 13).Value = 1782.57143
$ws.Cells.Item(107, 14).Value = -4212.5
$ws.Cells.Item(132, 8).Value = 11674770
$ws.Cells.Item(132, 9).Value = 11256383
$ws.Cells.Item(132, 10).Value = 12989698
$ws.Cells.Item(132, 11).Value = 33769149
$ws.Cells.Item(132, 12).Value = 38969094
$ws.Cells.Item(132, 13).Value = -33766619
$ws.Cells.Item(132, 14).Value = -38974154

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2572.7273
$ws.Cells.Item(40, 9).Value = 2025
$ws.Cells.Item(40, 10).Value = 2885.7144
$ws.Cells.Item(40, 11).Value = 2025
$ws.Cells.Item(40, 12).Value = 2885.7144
$ws.Cells.Item(40, 13).Value = -1889
$ws.Cells.Item(40, 14).Value = -3157.7144
$ws.Cells.Item(61, 8).Value = 1544.5454
$ws.Cells.Item(61, 9).Value = 1008.6667
$ws.Cells.Item(61, 10).Value = 2692.8572
$ws.Cells.Item(61, 11).Value = 1008.6667
$ws.Cells.Item(61, 12).Value = 2692.8572
$ws.Cells.Item(61, 13).Value = -806.6667
$ws.Cells.Item(61, 14).Value = -3096.8572
$ws.Cells.Item(113, 8).Value = 1544.5454
$ws.Cells.Item(113, 9).Value = 1008.6667
$ws.Cells.Item(113, 10).Value = 2692.8572
$ws.Cells.Item(113, 11).Value = 1008.6667
$ws.Cells.Item(113, 12).Value = 2692.8572
$ws.Cells.Item(113, 13).Value = 1161.3333
$ws.Cells.Item(113, 14).Value = -7032.8572
$ws.Cells.Item(132, 8).Value = 7946814
$ws.Cells.Item(132, 9).Value = 15886328
$ws.Cells.Item(132, 10).Value = 7299.778
$ws.Cells.Item(132, 11).Value = 47658984
$ws.Cells.Item(132, 12).Value = 21899.334
$ws.Cells.Item(132, 13).Value = -47656454
$ws.Cells.Item(132, 14).Value = -26959.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 450
$ws.Cells.Item(113, 9).Value = 250
$ws.Cells.Item(113, 11).Value = 750
$ws.Cells.Item(113, 13).Value = 1420
$ws.Cells.Item(136, 8).Value = 8473.833000000001
$ws.Cells.Item(136, 9).Value = 5129.0386
$ws.Cells.Item(136, 10).Value = 30215
$ws.Cells.Item(136, 11).Value = 15387.1158
$ws.Cells.Item(136, 12).Value = 90645
$ws.Cells.Item(136, 13).Value = -12837.1158
$ws.Cells.Item(136, 14).Value = -95745
